$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.405.95'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '3.503.04'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').Value = '3.501.73'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +2.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.01%  '
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.71%  '
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '4.098.65'
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '3.503.35'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '67.376.48'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.57'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '448.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.634'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').Value = '3.645.77'
$ws.Range('E25').Value = '  -0.50%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.86'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.11'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.66'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('E36').Value = '  +1.52%  '
$ws.Range('D37').Value = '3.498.18'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.08'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  +4.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0899'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.59%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '173.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.24%  '
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '30.33'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.00%  '
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E48').Value = '  +3.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.69'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.254'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.82%  '
